# Updates cryptos list data (prices, 1h volume %, and two swapped coin rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value renders as a plain number (e.g. "92.14") must be
# forced to Text format first, so Excel stores the literal string instead of
# silently converting it to a numeric value (these prices intentionally keep
# their original text formatting, matching the rest of the sheet).
$textFormatCells = @('D5', 'D6', 'D7', 'D10', 'D11', 'D15', 'D17', 'D21', 'D22', 'D23', 'D27', 'D28', 'D29', 'D30', 'D31', 'D34', 'D36', 'D38', 'D39', 'D44', 'D45', 'D47', 'D49', 'D50', 'D51')
foreach ($cell in $textFormatCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$changes = @{
    'D2' = '41.614.83'
    'E2' = '  +0.21%  '
    'D3' = '2.473.02'
    'E3' = '  -0.15%  '
    'E4' = '  +0.16%  '
    'D5' = '317.75'
    'E5' = '  +1.48%  '
    'D6' = '92.14'
    'E6' = '  -0.45%  '
    'D7' = '0.552'
    'E7' = '  +1.46%  '
    'E8' = '  +0.06%  '
    'E9' = '  +1.75%  '
    'D10' = '0.0862'
    'E10' = '  +10.17%  '
    'D11' = '33.05'
    'E11' = '  +1.25%  '
    'E12' = '  +0.68%  '
    'D13' = '2.854.56'
    'E13' = '  -0.14%  '
    'E14' = '  +0.77%  '
    'D15' = '15.63'
    'E15' = '  -3.51%  '
    'D16' = '2.472.21'
    'E16' = '  +0.09%  '
    'D17' = '0.788'
    'E17' = '  +2.83%  '
    'D18' = '41.571.70'
    'E18' = '  +0.11%  '
    'E19' = '  +0.63%  '
    'E20' = '  +1.20%  '
    'D21' = '71.27'
    'E21' = '  -0.89%  '
    'D22' = '11.33'
    'E22' = '  +0.68%  '
    'D23' = '240.01'
    'E23' = '  +1.48%  '
    'E24' = '  +1.25%  '
    'E25' = '  +1.77%  '
    'E26' = '  +0.02%  '
    'D27' = '24.69'
    'E27' = '  -0.44%  '
    'D28' = '2.28'
    'E28' = '  +3.59%  '
    'D29' = '9.85'
    'E29' = '  +2.41%  '
    'D30' = '36.21'
    'E30' = '  +1.19%  '
    'D31' = '160.72'
    'E31' = '  +1.47%  '
    'E32' = '  +1.49%  '
    'E33' = '  +0.03%  '
    'D34' = '0.0769'
    'E34' = '  +1.84%  '
    'E35' = '  +0.09%  '
    'D36' = '17.24'
    'E36' = '  +0.08%  '
    'E37' = '  +0.59%  '
    'B38' = 'ARBITRUM'
    'C38' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D38' = '1.83'
    'E38' = '  +1.11%  '
    'B39' = 'Stellar'
    'C39' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D39' = '0.116'
    'E39' = '  +1.51%  '
    'E40' = '  -3.49%  '
    'E41' = '  -3.20%  '
    'E42' = '  +2.20%  '
    'D43' = '1.991.28'
    'E43' = '  +0.50%  '
    'B44' = 'VeChain'
    'C44' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D44' = '0.0285'
    'E44' = '  +0.62%  '
    'B45' = 'EnergySwap'
    'C45' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D45' = '18.93'
    'E45' = '  -2.44%  '
    'E46' = '  +1.63%  '
    'D47' = '9.18'
    'E47' = '  +3.07%  '
    'D48' = '2.711.93'
    'E48' = '  -0.26%  '
    'D49' = '97.58'
    'E49' = '  +0.20%  '
    'B50' = 'BitcoinSV'
    'C50' = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
    'D50' = '73.49'
    'E50' = '  +1.68%  '
    'B51' = 'ordi'
    'C51' = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
    'D51' = '67.05'
    'E51' = '  -1.33%  '
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
